# Adjusting the NRG location
# Shift every timestamp in column A forward by exactly one day (the whole
# dataset now represents the next day), and update column B ("Actual
# Production (MW)") with the newly fetched production values for that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB = @(463,450,451,452,460,458,444,427,433,421,406,377,360,351,347,347,348,331,307,303,292,288,287,291,299,313,320,322,315,303,303,294,274,253,235,216,196,192,210,233,281,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $newB.Length; $i++) {
    $row = $i + 2
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $aCell.Value2 + 1
    $ws.Cells.Item($row, 2).Value = $newB[$i]
}
